$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.530.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.957.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.59%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.28%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.16"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.00%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.77%  "

# Row 12
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.99%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.843"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.239.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.56"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.74%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.957.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.92%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.397.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.08"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0852"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.54%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.55%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.38%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.83%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.84%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.24"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.13%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +19.63%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.78"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0611"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.18%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.32%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.77%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -13.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0966"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.88%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0209"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.43%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.78"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.79%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.359.06"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.74%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.19"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.53%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.84"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.134.67"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.53%  "
